# Workbook / worksheet handles
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# xlPasteFormats
$xlPasteFormats = -4122

# ---------------------------------------------------------------
# 1. Stamp formatting onto the new rows (52-65) by copying from
#    existing template rows that already carry the right styles.
#    Row 51 = "normal" day row (s=2,15,15,15,16)
#    Row 42 = "Away"   day row (s=2,18,19,19,20)
# ---------------------------------------------------------------
$ws.Range("A51:E51").Copy()
$ws.Range("A52:E52").PasteSpecial($xlPasteFormats)
$ws.Range("A53:E53").PasteSpecial($xlPasteFormats)
$ws.Range("A54:E54").PasteSpecial($xlPasteFormats)
$ws.Range("A55:E55").PasteSpecial($xlPasteFormats)
$ws.Range("A59:E59").PasteSpecial($xlPasteFormats)
$ws.Range("A60:E60").PasteSpecial($xlPasteFormats)
$ws.Range("A61:E61").PasteSpecial($xlPasteFormats)
$ws.Range("A62:E62").PasteSpecial($xlPasteFormats)
$ws.Range("A63:E63").PasteSpecial($xlPasteFormats)
$ws.Range("A64:E64").PasteSpecial($xlPasteFormats)
$ws.Range("A65:E65").PasteSpecial($xlPasteFormats)

$ws.Range("A42:E42").Copy()
$ws.Range("A56:E56").PasteSpecial($xlPasteFormats)
$ws.Range("A57:E57").PasteSpecial($xlPasteFormats)
$ws.Range("A58:E58").PasteSpecial($xlPasteFormats)

# Row 53's B/C/D cells use style 16 (same look but no vertical-center)
# instead of style 15 - pick that style up from an existing cell (C3)
# that already carries it.
$ws.Range("C3").Copy()
$ws.Range("B53").PasteSpecial($xlPasteFormats)
$ws.Range("C53").PasteSpecial($xlPasteFormats)
$ws.Range("D53").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------
# 2. Fill in the date serials (column A) for the new rows.
# ---------------------------------------------------------------
$ws.Range("A52").Value = 41359
$ws.Range("A53").Value = 41360
$ws.Range("A54").Value = 41361
$ws.Range("A55").Value = 41362
$ws.Range("A56").Value = 41363
$ws.Range("A57").Value = 41364
$ws.Range("A58").Value = 41365
$ws.Range("A59").Value = 41366
$ws.Range("A60").Value = 41367
$ws.Range("A61").Value = 41368
$ws.Range("A62").Value = 41369
$ws.Range("A63").Value = 41370
$ws.Range("A64").Value = 41371
$ws.Range("A65").Value = 41372

# ---------------------------------------------------------------
# 3. Fill in the workload readings (columns B-E).
# ---------------------------------------------------------------
$ws.Range("B52").Value = "0.5H"
$ws.Range("C52").Value = "0H"
$ws.Range("D52").Value = "0H"
$ws.Range("E52").Value = "-"

$ws.Range("B53").Value = "0.5H"
$ws.Range("C53").Value = "0H"
$ws.Range("D53").Value = "0H"
$ws.Range("E53").Value = "-"

$ws.Range("B54").Value = "0.5H"
$ws.Range("C54").Value = "0H"
$ws.Range("D54").Value = "0H"
$ws.Range("E54").Value = "-"

$ws.Range("B55").Value = "1H"
$ws.Range("C55").Value = "0H"
$ws.Range("D55").Value = "0H"
$ws.Range("E55").Value = "-"

$ws.Range("B56").Value = "Away"
$ws.Range("B57").Value = "Away"
$ws.Range("B58").Value = "Away"

$ws.Range("B59").Value = "0.5H"
$ws.Range("C59").Value = "0H"
$ws.Range("D59").Value = "0H"
$ws.Range("E59").Value = "-"

$ws.Range("B60").Value = "1H"
$ws.Range("C60").Value = "0H"
$ws.Range("D60").Value = "0H"
$ws.Range("E60").Value = "-"

$ws.Range("B61").Value = "0.5H"
$ws.Range("C61").Value = "0H"
$ws.Range("D61").Value = "0H"
$ws.Range("E61").Value = "-"

$ws.Range("B62").Value = "0.5H"
$ws.Range("C62").Value = "0H"
$ws.Range("D62").Value = "0H"
$ws.Range("E62").Value = "-"

$ws.Range("B63").Value = "1H"
$ws.Range("C63").Value = "0H"
$ws.Range("D63").Value = "0H"
$ws.Range("E63").Value = "-"

$ws.Range("B64").Value = "1.5H"
$ws.Range("C64").Value = "0H"
$ws.Range("D64").Value = "0H"
$ws.Range("E64").Value = "-"

$ws.Range("E65").Value = "-"

# ---------------------------------------------------------------
# 4. Merge the "Away" rows' B:E cells, matching existing rows such
#    as B42:E42.
# ---------------------------------------------------------------
$ws.Range("B56:E56").Merge()
$ws.Range("B57:E57").Merge()
$ws.Range("B58:E58").Merge()

# ---------------------------------------------------------------
# 5. Update the visible selection to the new bottom-right-most cell.
# ---------------------------------------------------------------
$ws.Range("B65").Select()
